$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.556.55"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.909.91"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.02"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.29"
$ws.Range("E6").Value = "  -6.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.48"
$ws.Range("E10").Value = "  -4.83%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.86"
$ws.Range("E13").Value = "  -5.36%  "
$ws.Range("D14").Value = "3.383.73"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "2.922.13"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.961"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "51.552.09"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.44"
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.66"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.40"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  -3.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("E26").Value = "  +8.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.170"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "26.37"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.103"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.16"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.09"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.30"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.72"
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0422"
$ws.Range("E37").Value = "  -6.30%  "
$ws.Range("E38").Value = "  -7.27%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.54"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.93"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "2.084.09"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  -8.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.236"
$ws.Range("E49").Value = "  -5.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0335"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.887"
$ws.Range("E51").Value = "  -6.03%  "
